# Update the "Automation" section's first task row (row 24 on the Tasks sheet).
# The task text was revised, Story Points / (hidden) estimate columns were
# filled in, and a new Description note was added explaining the relationship
# consideration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("B24").Value = "Automatically update Last Odometer reading, -date of Vehicle records based on last created Vehicle Renting  record"
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = "Consider add additional relationship and you may need to use validation rule to ensure correct vehicle record will be updated"

# Move the active selection to match the author's saved cursor position.
$ws.Range("C25").Select()
